# Week 2 practice-problem solutions: refresh the pasted screenshot's
# stored display size (the picture was resized in the Word UI) and
# switch the page margins over to the narrower "Normal (0.5")"-ish preset.

$d = $word.ActiveDocument

# --- 1. Resize the inline picture -----------------------------------
# Grab the (only) inline picture in the document.
$shp = $d.InlineShapes.Item(1)
$shpRange = $shp.Range

# Build the replacement run as real WordprocessingML (same run, but with
# the new drawing extents/effectExtent/editId and an explicit
# <w:rPr><w:noProof/></w:rPr> on the run, as Word stamps on picture runs)
# and feed it through the Range.InsertXML COM entry point - the
# Word-native way of dropping raw OOXML into the story.
$newRunXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
  'xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" ' +
  'xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" ' +
  'xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" ' +
  'xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" ' +
  'xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing">' +
  '<w:body>' +
  '<w:p>' +
  '<w:r w:rsidRPr="009F039C">' +
  '<w:rPr><w:noProof/></w:rPr>' +
  '<w:drawing>' +
  '<wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="1924DBB2" wp14:editId="6163E9E6">' +
  '<wp:extent cx="6656749" cy="1747157"/>' +
  '<wp:effectExtent l="0" t="0" r="0" b="5715"/>' +
  '<wp:docPr id="1724334495" name="Picture 1"/>' +
  '<wp:cNvGraphicFramePr>' +
  '<a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/>' +
  '</wp:cNvGraphicFramePr>' +
  '<a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main">' +
  '<a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">' +
  '<pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture">' +
  '<pic:nvPicPr>' +
  '<pic:cNvPr id="1724334495" name=""/>' +
  '<pic:cNvPicPr/>' +
  '</pic:nvPicPr>' +
  '<pic:blipFill>' +
  '<a:blip r:embed="rId4"/>' +
  '<a:stretch><a:fillRect/></a:stretch>' +
  '</pic:blipFill>' +
  '<pic:spPr>' +
  '<a:xfrm><a:off x="0" y="0"/><a:ext cx="6688621" cy="1755522"/></a:xfrm>' +
  '<a:prstGeom prst="rect"><a:avLst/></a:prstGeom>' +
  '</pic:spPr>' +
  '</pic:pic>' +
  '</a:graphicData>' +
  '</a:graphic>' +
  '</wp:inline>' +
  '</w:drawing>' +
  '</w:r>' +
  '</w:p>' +
  '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData>' +
  '</pkg:part>' +
  '</pkg:package>'

# Remove the old drawing run, then drop the rebuilt run in its place
# (re-using the same Range object keeps the insertion inside the
# existing paragraph instead of splitting a new one off).
$shpRange.Delete()
$shpRange.InsertXML($newRunXml)

# --- 2. Narrower page margins ----------------------------------------
$ps = $d.PageSetup
$ps.TopMargin = 36
$ps.BottomMargin = 36
$ps.LeftMargin = 36
$ps.RightMargin = 36

Write-Output "Picture resized and margins updated."
